# TC06_Canine_Filter_Breed-BelgMalin.xlsx edit
# Commit: "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The CasesTab query (B2) and FilesTab query (B4) had been swapped in the
# shared-string table; on top of that the CasesTab query itself had a stray
# trailing `co.cohort_description` (Cohort) line that does not belong to the
# query and needs to be removed. This script restores the FilesTab query to
# B4 and writes the corrected (Cohort-less) CasesTab query back into B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# B2 : CasesTab query - corrected text (no more "Cohort" column)
# ---------------------------------------------------------------------
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Belgian Malinois']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# ---------------------------------------------------------------------
# B4 : FilesTab query - unchanged text
# ---------------------------------------------------------------------
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Belgian Malinois']  
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Here-strings keep a trailing newline - trim it so the text matches exactly.
$casesQuery = $casesQuery.TrimEnd("`r", "`n")
$filesQuery = $filesQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $casesQuery
$ws.Range("B4").Value = $filesQuery

# ---------------------------------------------------------------------
# Row heights follow from the new wrapped-text line counts (17 lines for
# the Cases query, 16 for the unchanged Samples query, 15 for the Files
# query), combined with Excel's slightly reduced per-line height.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# ---------------------------------------------------------------------
# View state: zoomed in to 115%, scrolled/selected on row 2 now.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
